$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 with the new error message text (replacing "R2")
$ws.Range("B2").Value = "java.net.ConnectException: Connection timed out: connect 8888"

# Remove row 3 entirely (it held the old standalone "java.net.ConnectException..." text in A3)
$ws.Rows("3").Delete()

# Update the active selection to C3, matching the edited workbook view state
$ws.Range("C3").Select()
